$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.255.00"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "3.499.09"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "4.096.92"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").Value = "3.500.39"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.96%  "
$ws.Range("D17").Value = "64.347.75"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.572"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "3.639.57"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -4.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").Value = "3.520.99"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.68%  "
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0783"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.31%  "
$ws.Range("D48").Value = "2.463.38"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("E51").Value = "  -1.97%  "
